$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31: num_customers 46 -> 47 (cohort_size stays 2312), retention_rate recalculated
$ws.Range("C31").Value2 = 47
$ws.Range("E31").Value2 = $ws.Range("C31").Value2 / $ws.Range("D31").Value2

# Row 37: num_customers 832 -> 834, cohort_size 832 -> 834 (retention_rate stays 1)
$ws.Range("C37").Value2 = 834
$ws.Range("D37").Value2 = 834
